# Converts an "RRGGBB" hex string into the BGR-packed Long that the
# PowerPoint COM `RGB` property (ThemeColor.RGB / RGB()) expects
# (0x00BBGGRR).
function Convert-HexToRgbLong {
    param([string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# --- 1. Re-colour the deck's theme (slide master) from the "Integral"
#        palette to the stock "Office Theme" palette -------------------
# Order matches ThemeColorScheme.Colors(1..12):
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $masterScheme.Count; $i++) {
    $masterScheme.Item($i).RGB = Convert-HexToRgbLong $officeColors[$i - 1]
}

# --- 2. Update the cash-flow table on slide 16 to use the new table
#        style referenced by the edit -----------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{FCBE5D27-2C39-4779-9067-7C67C544B6AE}")
    }
}
